$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

foreach ($ws in @($ws1, $ws4)) {
    $ws.Range("F2").Value = 2342
    $ws.Range("G2").Value = "不可售"
    $ws.Range("F3").Value = 1813
    $ws.Range("F4").Value = 350
    $ws.Range("F5").Value = 1118
    $ws.Range("F6").Value = 1019
    $ws.Range("F8").Value = 5909
    $ws.Range("F9").Value = 94
}
